$d = $word.ActiveDocument

# The log entry paragraph ("Jag arbetar i grupp...") currently carries the
# "FirstParagraph" style (it is the first entry in the logbook). The edit adds
# a new first line with the date, and the existing text becomes an ordinary
# "BodyText" paragraph right below it.
$logText = "Jag arbetar i grupp med Alexander, och som instruerade ska vi bygga paperstorn."
$dateText = "[24/8/20]"

# Find the paragraph holding the log entry text.
$logPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    $candidateText = $candidate.Range.Text.TrimEnd([char]13, [char]7)
    if ($candidateText -eq $logText) {
        $logPara = $candidate
        break
    }
}

if ($logPara -eq $null) {
    throw "Could not find the log-entry paragraph to split."
}

# The paragraph right after the log entry already carries the "BodyText"
# style. Re-typing the log text at the very start of that paragraph (instead
# of changing the original paragraph's style) lets the new paragraph pick up
# "BodyText" formatting naturally, without disturbing run formatting further
# down (e.g. the bold "Trianglar" run) in that paragraph.
$nextPara = $logPara.Next()
$insertionPoint = $d.Range($nextPara.Range.Start, $nextPara.Range.Start)
$insertionPoint.InsertBefore($logText)

# Split right after the text we just inserted so it becomes its own
# "BodyText" paragraph, separate from the original body text that follows.
$splitAt = $insertionPoint.Start + $logText.Length
$splitRange = $d.Range($splitAt, $splitAt)
$splitRange.InsertParagraphBefore()

# The original paragraph (still styled "FirstParagraph") now just needs its
# text swapped out for the date stamp.
$logPara.Range.Text = $dateText
